$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Query 1 SparkSQL")

# Update Hot restart execution-time samples (column A, rows 3-101)
$aValues = @(935,867,858,825,822,794,772,790,775,791,758,791,773,760,779,797,782,775,780,761,776,782,770,769,766,767,774,769,761,781,769,778,774,771,781,763,758,782,761,756,773,768,755,782,767,774,758,768,772,772,766,772,770,774,760,766,753,760,758,764,769,750,774,784,771,779,773,765,759,756,761,770,762,760,772,766,767,764,764,784,774,758,765,772,772,765,771,770,828,775,773,767,756,766,761,782,768,786,768)
for ($i = 0; $i -lt $aValues.Length; $i++) {
    $ws.Cells.Item($i + 3, 1).Value = $aValues[$i]
}

# Update Cold restart execution-time samples (column E, rows 3-8)
$eValues = @(6452,6177,6125,6377,6140,6103)
for ($i = 0; $i -lt $eValues.Length; $i++) {
    $ws.Cells.Item($i + 3, 5).Value = $eValues[$i]
}

# Make this sheet the active tab and park the selection on E8
[void]$ws.Activate()
[void]$ws.Range("E8").Select()
